# Edit script: rename headers and add PO Forecast sheet with forecast data

$wb = $excel.ActiveWorkbook

# 1. Rename header on "Weekly Quantity" sheet
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header on "Monthly Trend" sheet
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add new "PO Forecast" worksheet after "Monthly Trend".
#    Copy the "Monthly Trend" sheet so the page setup / outline props match
#    the rest of the workbook, then clear it out and rename it.
$wsMonthly.Copy($null, $wsMonthly)
$wsForecast = $wb.Worksheets.Item("Monthly Trend (2)")
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.Clear()

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header style from Monthly Trend A1:B1 (bold, bordered, centered) to A1:D1
$wsMonthly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$data = @(
    @(44934.99999999999, 4,  3.999973161566078,  3.999973161661693),
    @(44941.99999999999, 14, 13.99997716231689,  13.99997716241489),
    @(44948.99999999999, 24, 23.99998104949662,  23.9999812769637),
    @(44955.99999999999, 34, 33.99998478948518,  33.9999855423474),
    @(44962.99999999999, 44, 43.99998843550167,  43.99998984535112),
    @(44969.99999999999, 54, 53.99999201671665,  53.99999421331075),
    @(44976.99999999999, 64, 63.99999557160424,  63.99999864698897),
    @(44983.99999999999, 74, 73.9999990274157,   74.00000317159515),
    @(44990.99999999999, 84, 84.00000242446531,  84.0000077171219),
    @(44997.99999999999, 94, 94.00000594194492,  94.00001240408955)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply date-style (matching column A style from Monthly Trend) to the ds column data cells
$wsMonthly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()
